$d = $word.ActiveDocument
$wdColorRed = 255

function Set-ParagraphRed($para) {
    $para.Range.Font.Color = $wdColorRed
}

# 1) "Consulta de saldo mora. Mora de un socio por día."
#    currently colored accent2 (C0504D/themeColor accent2) -> plain red FF0000
$rng = $d.Content
$found = $rng.Find.Execute("Consulta de saldo mora. Mora de un socio por día.")
if ($found) {
    Set-ParagraphRed $rng.Paragraphs(1)
}

# 2) "Líneas y tipos de operaciones asociadas" -> add red color
#    3) "Contragarantías" is the very next paragraph (sub-bullet); the word
#    "contragarantías" also appears earlier in the document inside a plain
#    sentence, so rather than searching for it directly (ambiguous), anchor
#    on the unique heading above it and move to the following paragraph.
$rng = $d.Content
$found = $rng.Find.Execute("Líneas y tipos de operaciones asociadas")
if ($found) {
    $lineasPara = $rng.Paragraphs(1)
    Set-ParagraphRed $lineasPara
    $contragarantiasPara = $lineasPara.Next()
    Set-ParagraphRed $contragarantiasPara
}

# 4) "Desembolsos y recuperos" and 5) "Administración y consultas"
#    These two paragraphs are adjacent; "Administración y consultas" is unique
#    in the document, so find it and then step back to the previous paragraph
#    to reach "Desembolsos y recuperos".
$rng = $d.Content
$found = $rng.Find.Execute("Administración y consultas")
if ($found) {
    $adminPara = $rng.Paragraphs(1)
    $prevPara = $adminPara.Previous()
    Set-ParagraphRed $prevPara
    Set-ParagraphRed $adminPara
}
